$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 66: latest EUR->ARS quote (2025-10-08T21:27:00Z).
# Force column A to Text format before assignment so the date-like
# string "2025-10-08" is stored as text, matching the other rows in
# this column instead of being auto-converted to a date serial number.
$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = "2025-10-08"
$ws.Range("B66").Value = "21:27:00"
$ws.Range("C66").Value = "1.00 EUR = 1,770.2348"
